$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("upload")

# Mejoras al netobook 5: corregir el texto del stub period y actualizar
# la selección activa a la columna F (tenor/stub_period) en vez de C.
$ws.Cells.Replace("SHORTBACK", "SHORTFRONT")

$ws.Activate()
$ws.Range("F2:F34").Select()
